$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.28"
$ws.Range("E2").Value = "'-0.69%"
$ws.Range("D3").Value = "'27.21"
$ws.Range("E3").Value = "'2.94%"
$ws.Range("D4").Value = "'5.116"
$ws.Range("E4").Value = "'0.79%"
$ws.Range("D5").Value = "'0.05714"
$ws.Range("E5").Value = "'2.10%"
$ws.Range("D6").Value = "'6.492"
$ws.Range("E6").Value = "'0.03%"
$ws.Range("D7").Value = "'0.8188"
$ws.Range("E7").Value = "'0.76%"
$ws.Range("D8").Value = "'0.8519"
$ws.Range("E8").Value = "'0.86%"
$ws.Range("B9").Value = "One"
$ws.Range("C9").Value = "https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("D9").Value = "'0.01015"
$ws.Range("E9").Value = "'1,603.48%"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1328"
$ws.Range("E10").Value = "'-0.85%"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.06934"
$ws.Range("E11").Value = "'-0.73%"
$ws.Range("B12").Value = "BitrueCoin"
$ws.Range("C12").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D12").Value = "'0.02827"
$ws.Range("E12").Value = "'-0.75%"
$ws.Range("B13").Value = "BitMartToken"
$ws.Range("C13").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D13").Value = "'0.09395"
$ws.Range("E13").Value = "'0.06%"
$ws.Range("B14").Value = "BitForexToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D14").Value = "'0.001522"
$ws.Range("E14").Value = "'0.42%"
$ws.Range("B15").Value = "CoinExToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D15").Value = "'0.04031"
$ws.Range("E15").Value = "'-13.41%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.006201"
$ws.Range("E16").Value = "'0.85%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.512"
$ws.Range("E17").Value = "'-2.54%"
$ws.Range("B18").Value = "GateToken"
$ws.Range("C18").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D18").Value = "'3.007"
$ws.Range("E18").Value = "'-0.25%"
$ws.Range("B19").Value = "BTSEToken"
$ws.Range("C19").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D19").Value = "'2.309"
$ws.Range("E19").Value = "'12.32%"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "'0.3165"
$ws.Range("E20").Value = "'0.29%"
$ws.Range("D21").Value = "'0.03220"
$ws.Range("E21").Value = "'0.68%"
$ws.Range("D22").Value = "'0.1303"
$ws.Range("E22").Value = "'0.47%"
$ws.Range("D23").Value = "'3.566"
$ws.Range("E23").Value = "'-5.15%"
$ws.Range("E24").Value = "'1.72%"
$ws.Range("D25").Value = "'0.001217"
$ws.Range("E25").Value = "'-2.04%"
$ws.Range("D26").Value = "'0.004475"
$ws.Range("E26").Value = "'-2.51%"
$ws.Range("D27").Value = "'0.00009891"
$ws.Range("E27").Value = "'3.08%"
$ws.Range("D28").Value = "'0.0001448"
$ws.Range("E28").Value = "'3.62%"
$ws.Range("D40").Value = "'0.03726"
$ws.Range("B41").Value = "KickToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
$ws.Range("D41").Value = "'0.005855"
$ws.Range("E41").Value = "'-5.07%"
$ws.Range("B42").Value = "BKEXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
$ws.Range("D42").Value = "'0.1059"
$ws.Range("E42").Value = "'0.07%"
$ws.Range("B43").Value = "CEJI"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbKjCVJCh+ceji-ceji"
$ws.Range("D43").Value = "'0.002298"
$ws.Range("E43").Value = "'-8.04%"
$ws.Range("D44").Value = "'0.009405"
$ws.Range("E44").Value = "'5.15%"
$ws.Range("D45").Value = "'0.00005150"
$ws.Range("E45").Value = "'-2.38%"
$ws.Range("E46").Value = "'-0.04%"
$ws.Range("E47").Value = "'-7.77%"
$ws.Range("D48").Value = "'0.002507"
$ws.Range("E48").Value = "'-4.45%"
$ws.Range("E49").Value = "'-0.04%"
$ws.Range("E50").Value = "'-0.04%"

# Clear the quote-prefix / number style introduced by forcing text values,
# so cells match plain (unstyled) string cells as in the source workbook.
$ws.Range("B2:E50").Style = "Normal"
